$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("A1").Value = "testUrls"
$ws.Range("B1").Value = "mainUrls"

# Move active selection to C3 (as in the diff)
$ws.Range("C3").Select()

$wb.Save()
